$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.767.51"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").Value = "3.388.56"
$ws.Range("E3").Value = "  -1.86%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.17%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.388.40"
$ws.Range("E8").Value = "  -1.90%  "

$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("E10").Value = "  -1.94%  "

$ws.Range("E11").Value = "  -1.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.398"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.82%  "

$ws.Range("D13").Value = "3.965.39"
$ws.Range("E13").Value = "  -1.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("E15").Value = "  +1.88%  "

$ws.Range("E16").Value = "  -2.06%  "

$ws.Range("D17").Value = "3.383.43"
$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("D18").Value = "60.858.55"
$ws.Range("E18").Value = "  -1.46%  "

$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.17%  "

$ws.Range("E21").Value = "  -5.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.88%  "

$ws.Range("E23").Value = "  -1.13%  "

$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("E26").Value = "  -5.49%  "

$ws.Range("D27").Value = "3.523.69"
$ws.Range("E27").Value = "  -1.91%  "

$ws.Range("E28").Value = "  -1.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.95%  "

$ws.Range("E31").Value = "  -2.43%  "

$ws.Range("E32").Value = "  -2.36%  "

$ws.Range("E33").Value = "  -3.25%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.69%  "

$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("E38").Value = "  -2.37%  "

$ws.Range("D39").Value = "3.418.39"
$ws.Range("E39").Value = "  -1.74%  "

$ws.Range("E40").Value = "  -4.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.782"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.16%  "

$ws.Range("E47").Value = "  -2.82%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.524.30"
$ws.Range("E48").Value = "  -2.02%  "

$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.65%  "
